# GI05MOAS GL484 Ingest CSV
# Created Ingest CSV for GL484 D2
#
# The workbook's Asset_Cal_Info sheet still carried the template's
# GI05MOAS-GL001 reference designators in column A ("Ref Des"). Update
# them to the GL484 deployment-2 instrument ref-des values, matching the
# already-populated Mooring Serial Number (484) / Deployment Number (2)
# columns on the same rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# FLORD (fluorometer) rows - calibration coefficient rows 2-5
$ws.Range("A2").Value = "GI05MOAS-GL484-01-FLORDM000"
$ws.Range("A3").Value = "GI05MOAS-GL484-01-FLORDM000"
$ws.Range("A4").Value = "GI05MOAS-GL484-01-FLORDM000"
$ws.Range("A5").Value = "GI05MOAS-GL484-01-FLORDM000"

# DOSTA (dissolved oxygen) row 7
$ws.Range("A7").Value = "GI05MOAS-GL484-02-DOSTAM000"

# CTDGV (CTD) row 9
$ws.Range("A9").Value = "GI05MOAS-GL484-04-CTDGVM000"

# ENG (engineering/glider) row 11
$ws.Range("A11").Value = "GI05MOAS-GL484-00-ENG000000"

# Restore the active sheet / selection state recorded in the saved file.
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Range("E28").Select()
$ws.Activate()
$ws.Range("B16").Select()
